# Training Dashboard progress update - 04 Nov 2025
#
# For rows 3-7:
#   - column H ("PERIOD TO EXPIRE"): decrement by 1 day
#   - column I ("LAST UPDATE"): bump the date text from 03-Nov-2025 to 04-Nov-2025
#
# The "LAST UPDATE" column stores its date as literal text (not a real Excel
# date serial). Assigning a date-like string straight to `.Value` makes Excel
# auto-convert it into a date serial (and reformats the cell), which would
# not match the source data's plain-text representation. To keep it as text
# with the original formatting untouched, we build the text via a formula in
# a scratch cell and paste its computed value back as a value-only paste,
# which Excel keeps as text instead of re-interpreting it as a date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$scratch = $ws.Range("Z1")

$rows = 3, 4, 5, 6, 7

foreach ($r in $rows) {
    $hCell = $ws.Cells.Item($r, 8)   # column H - PERIOD TO EXPIRE
    $hCell.Value = $hCell.Value() - 1

    $iCell = $ws.Cells.Item($r, 9)   # column I - LAST UPDATE
    $scratch.Formula = "=""04-Nov-2025"""
    $scratch.Copy()
    $iCell.PasteSpecial(-4163)       # xlPasteValues - keeps text as text
}

$scratch.Value = ""
$excel.CutCopyMode = 0
